$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "type index" lookup row (row 22), C22:T22, as shared-string text
$ws.Range("C22").Value = '"normal":0'
$ws.Range("D22").Value = '"flight":1'
$ws.Range("E22").Value = '"flying":2'
$ws.Range("F22").Value = '"poison":3'
$ws.Range("G22").Value = '"ground":4'
$ws.Range("H22").Value = '"rock":5'
$ws.Range("I22").Value = '"bug":6'
$ws.Range("J22").Value = '"ghost":7'
$ws.Range("K22").Value = '"steel":8'
$ws.Range("L22").Value = '"fire":9'
$ws.Range("M22").Value = '"water":10'
$ws.Range("N22").Value = '"grass":11'
$ws.Range("O22").Value = '"electric":12'
$ws.Range("P22").Value = '"psychic":13'
$ws.Range("Q22").Value = '"ice":14'
$ws.Range("R22").Value = '"dragon":15'
$ws.Range("S22").Value = '"dark":16'
$ws.Range("T22").Value = '"fairy":17'

# Update the selected cell shown when the workbook is reopened
$ws.Range("T20").Select()
